$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Switch to manual calculation before editing so that the existing shared
# formula in column B (SUM(C+D)) is not recalculated against the new text
# values being placed in D6/E6 - the source edit only changed the two
# target cells and left every other cached cell value untouched.
$excel.Calculation = -4135   # xlCalculationManual

# Row 6 ("4/15/2023"): replace the numeric entries in D6 and E6 with the
# literal text values "4:36" and "17:05" (these become new shared-string
# table entries rather than numbers/times).
$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "17:05"
